# Apply updated crypto price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as text in the source feed (they use
# '.' as a thousands separator in some rows), so force text format before
# assigning to prevent Excel from auto-converting them to numbers.
$dRanges = @("D2","D3","D5","D6","D7","D10","D12","D14","D15","D17","D18","D19","D21","D22","D23","D24","D27","D28","D29","D30","D32","D33","D34","D35","D39","D40","D41","D44","D45","D46","D48","D50","D51")
foreach ($addr in $dRanges) {
    $ws.Range($addr).NumberFormat = "@"
}

$values = @{
    'D2' = '44.233.75'
    'D3' = '2.254.29'
    'D5' = '308.15'
    'D6' = '98.87'
    'D7' = '0.575'
    'D10' = '35.58'
    'D12' = '7.33'
    'D14' = '2.598.36'
    'D15' = '2.277.54'
    'D17' = '13.84'
    'D18' = '44.107.93'
    'D19' = '12.92'
    'D21' = '6.34'
    'D22' = '65.50'
    'D23' = '243.13'
    'D24' = '2.95'
    'D27' = '10.13'
    'D28' = '2.14'
    'D29' = '36.75'
    'D30' = '6.19'
    'D32' = '3.58'
    'D33' = '157.18'
    'D34' = '0.0826'
    'D35' = '2.67'
    'D39' = '15.25'
    'D40' = '3.90'
    'D41' = '3.39'
    'D44' = '1.771.95'
    'D45' = '88.42'
    'D46' = '5.18'
    'D48' = '101.58'
    'D50' = '70.43'
    'D51' = '55.46'
    'E2' = '  +1.42%  '
    'E3' = '  +0.50%  '
    'E4' = '  +0.10%  '
    'E5' = '  -4.48%  '
    'E6' = '  -2.18%  '
    'E7' = '  -0.57%  '
    'E8' = '  +0.11%  '
    'E9' = '  -3.22%  '
    'E10' = '  -4.63%  '
    'E11' = '  -1.11%  '
    'E12' = '  -4.93%  '
    'E13' = '  -1.63%  '
    'E14' = '  +0.41%  '
    'E15' = '  +1.34%  '
    'E16' = '  -1.75%  '
    'E17' = '  -2.49%  '
    'E18' = '  +1.25%  '
    'E19' = '  -5.12%  '
    'E20' = '  -0.80%  '
    'E21' = '  -3.10%  '
    'E24' = '  -6.58%  '
    'E25' = '  -8.47%  '
    'E26' = '  +0.18%  '
    'E27' = '  +0.49%  '
    'E28' = '  -1.88%  '
    'E29' = '  -0.13%  '
    'E30' = '  -1.73%  '
    'E31' = '  +0.10%  '
    'E32' = '  +13.98%  '
    'E33' = '  -1.93%  '
    'E34' = '  -3.35%  '
    'E35' = '  -0.30%  '
    'E36' = '  -0.43%  '
    'E37' = '  -4.25%  '
    'E38' = '  -3.79%  '
    'E39' = '  -2.25%  '
    'E40' = '  -8.19%  '
    'E41' = '  -10.30%  '
    'E42' = '  -3.17%  '
    'E43' = '  +0.07%  '
    'E44' = '  -1.48%  '
    'E45' = '  +7.55%  '
    'E46' = '  -0.39%  '
    'E47' = '  -2.99%  '
    'E48' = '  -1.42%  '
    'E49' = '  -1.83%  '
    'E50' = '  -5.21%  '
    'E51' = '  -5.21%  '
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
